$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab16")

# --- Fix mangled (mojibake) accented characters in the Regional Economic
# Communities footnote (PALOP / MERCOSUR definitions) ---
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = `"Community of Sahel-Saharan States`";COMESA = `"Common Market for Eastern and Southern Africa`";EAC = `"East African Community`";ECCAS = `"Economic Community of Central African States`";ECOWAS = `"Economic Community of West African States`";IGAD = `"Intergovernmental Authority on Development`";SADC = `"Southern African Development Community`";UMA = `"Arab Maghreb Union`";PALOP = `"Países Africanos de Língua Oficial Portuguesa`";ASEAN = `"Association of Southeast Asian Nations`";MERCOSUR = `"Mercado Común del Sur`".EU27 = `"European Union (27 members)`".OECD = `"Organisation for Economic Co-operation and Development`"."

# --- Minor recalculated value tweak ---
$ws.Range("G68").Value = 72.002804312094995

# --- Updated data for row 97: Africa, Fragile States ---
$ws.Range("C97").Value = 15.287622542340101
$ws.Range("D97").Value = 10.0156023030087
$ws.Range("E97").Value = 0.50664809386873999
$ws.Range("F97").Value = 19.1167573542069
$ws.Range("G97").Value = 51.512191835951199
$ws.Range("H97").Value = 177720.91648927901
$ws.Range("I97").Value = 117565.86665001699
$ws.Range("J97").Value = 6229.7809980038901
$ws.Range("K97").Value = 227426.29069014301
$ws.Range("L97").Value = 599259.95106765197

# --- Updated data for row 98: ROW, Fragile States ---
$ws.Range("C98").Value = 10.8971261517
$ws.Range("D98").Value = 7.8425138939590999
$ws.Range("E98").Value = 0.17124113348757999
$ws.Range("F98").Value = 15.594154452669899
$ws.Range("G98").Value = 56.889300078257001
$ws.Range("H98").Value = 163775.14051867899
$ws.Range("I98").Value = 117300.959905896
$ws.Range("J98").Value = 3037.3478936900601
$ws.Range("K98").Value = 229150.68862263701
$ws.Range("L98").Value = 811990.48959234694
